$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.277.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.930.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.56%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7516"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.07%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'249.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9983"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'Cardano"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'0.3227"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.23%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'Solana"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'28.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.41%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -3.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08011"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.35%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.930.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.394"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.78%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'94.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.65%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.282.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.23%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'253.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.20%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008049"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.44%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.767"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.188.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.26%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9985"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.9992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.840"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.92%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.587"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.74%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'164.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.85%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.322"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.49%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Stellar"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.1343"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.61%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'EthereumClassic"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'19.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.74%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.353"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.41%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D32").Value = "'4.426"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.154"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.303"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.33%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.05116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.88%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7484"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.764"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01973"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.07%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.798"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'78.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4504"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.38%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.9986"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'101.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.01%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'9.819"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.44%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Maker"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.003.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +13.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.01%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'37.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06065"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.17%  "
$ws.Range("E51").Style = "Normal"
